$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new log entry as a plain text string (matching the existing
# inline-string rows) into the next empty row of column A.
$ws.Cells.Item(24, 1).Value = "2025-10-26 23:49:06"
